# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as TEXT (e.g. "57.432.57" uses dots as thousands
# separators). Values that happen to look like a plain decimal number
# (e.g. "502.48") would otherwise be auto-converted to a numeric value by
# Excel, so those are written with a leading apostrophe to force text,
# exactly like typing `'502.48` into a cell in Excel.

$ws.Range('D2').Value = '57.432.57'
$ws.Range('E2').Value = '  -1.27%  '
$ws.Range('D3').Value = '2.426.83'
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''502.48'
$ws.Range('E5').Value = '  -3.45%  '
$ws.Range('D6').Value = '''127.71'
$ws.Range('E6').Value = '  -3.50%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('D9').Value = '2.438.90'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').Value = '''0.0944'
$ws.Range('E11').Value = '  -5.02%  '
$ws.Range('E12').Value = '  -4.13%  '
$ws.Range('D13').Value = '''0.326'
$ws.Range('E13').Value = '  -4.78%  '
$ws.Range('D14').Value = '2.862.43'
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('D15').Value = '57.366.55'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '''21.60'
$ws.Range('E16').Value = '  -2.19%  '
$ws.Range('E17').Value = '  -3.62%  '
$ws.Range('D18').Value = '2.436.59'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').Value = '''10.37'
$ws.Range('E19').Value = '  -4.60%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''313.13'
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = '''4.07'
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('D23').Value = '''5.66'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('D25').Value = '''0.403'
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = '''0.158'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').Value = '''7.14'
$ws.Range('E28').Value = '  -3.64%  '
$ws.Range('D29').Value = '''169.23'
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('D30').Value = '0.0₃0716'
$ws.Range('E30').Value = '  -4.61%  '
$ws.Range('D31').Value = '''6.15'
$ws.Range('E31').Value = '  -3.48%  '
$ws.Range('E32').Value = '  -3.53%  '
$ws.Range('D33').Value = '''1.11'
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '''0.998'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '''17.63'
$ws.Range('E36').Value = '  -2.94%  '
$ws.Range('D37').Value = '''1.25'
$ws.Range('E37').Value = '  -5.83%  '
$ws.Range('D38').Value = '''3.87'
$ws.Range('E38').Value = '  -3.01%  '
$ws.Range('D39').Value = '''36.43'
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('D41').Value = '''0.744'
$ws.Range('E41').Value = '  -6.36%  '
$ws.Range('D42').Value = '''268.77'
$ws.Range('E42').Value = '  -2.75%  '
$ws.Range('D43').Value = '''3.33'
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('D44').Value = '''4.82'
$ws.Range('E44').Value = '  -4.70%  '
$ws.Range('D45').Value = '''0.577'
$ws.Range('E45').Value = '  -3.28%  '
$ws.Range('D46').Value = '''0.0906'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').Value = '''118.41'
$ws.Range('E47').Value = '  -6.24%  '
$ws.Range('D48').Value = '''0.0482'
$ws.Range('E48').Value = '  -2.24%  '
$ws.Range('D49').Value = '''17.06'
$ws.Range('E49').Value = '  -4.62%  '
$ws.Range('D50').Value = '''0.0207'
$ws.Range('E50').Value = '  -3.88%  '
$ws.Range('D51').Value = '''16.41'
$ws.Range('E51').Value = '  -4.42%  '
